$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: omicwas
$ws.Range("A8").Value = "omicwas"
$ws.Range("B8").Value = 33752591
$ws.Range("C8").Value = "Nonlinear ridge regression improves cell-type-specific differential expression analysis"
$ws.Range("D8").Value = "R package"
$ws.Range("E8").Value = "Fumihiko Takeuchi"
$ws.Range("F8").Value = "Fumihiko Takeuchi [fumihiko@takeuchi.name]"
$ws.Range("G8").Value = "see running-notes.md"

# B8 previously held a leftover (now-removed) hyperlink style; reproduce that
# exact "used to be a hyperlink" look (no underline / no explicit color) by
# copying the formatting from an existing cell that already has it (G2),
# rather than resynthesizing font attributes by hand.
$ws.Range("G2").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# Row 9: TOAST
$ws.Range("A9").Value = "TOAST"
$ws.Range("B9").Value = 31484546
$ws.Range("C9").Value = "TOAST: improving reference-free cell composition estimation by cross-cell type differential analysis"
$ws.Range("D9").Value = "R package"
$ws.Range("E9").Value = "Ziyi Li"
$ws.Range("F9").Value = "Hao Wu [hao.wu@emory.edu]"
$ws.Range("G9").Value = "see running-notes.md"

# Update the cursor/selection to match where the author ended up (E10).
$ws.Range("E10").Select()
